# Update "想去人数" (want-to-go count) figures in the F column across sheets,
# matching the refreshed data pull reflected in the commit diff.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 10419
$ws1.Range("F6").Value = 613
$ws1.Range("F13").Value = 538
$ws1.Range("F20").Value = 415
$ws1.Range("F21").Value = 415
$ws1.Range("F25").Value = 1080
$ws1.Range("F31").Value = 259
$ws1.Range("F34").Value = 674

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 73
$ws2.Range("F21").Value = 2264
$ws2.Range("F46").Value = 87

# 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F10").Value = 413
$ws3.Range("F11").Value = 370

# 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 10419
$ws4.Range("F10").Value = 370
$ws4.Range("F19").Value = 415
$ws4.Range("F20").Value = 415
$ws4.Range("F24").Value = 2264
$ws4.Range("F25").Value = 2264
$ws4.Range("F27").Value = 1080
$ws4.Range("F34").Value = 674
$ws4.Range("F51").Value = 87
